# Refresh the cryptos price/volume snapshot and restore the Avalanche/Polygon row order
# (values pulled verbatim from the upstream coinranking.com feed, including its
# "dotted" thousands-separator price strings e.g. "37.306.05").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.306.05'
$ws.Cells.Item(2, 5).Value = '  +0.15%  '

$ws.Cells.Item(3, 4).Value = '2.008.04'
$ws.Cells.Item(3, 5).Value = '  -0.29%  '

$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$ws.Cells.Item(5, 4).Value = "'258.00"
$ws.Cells.Item(5, 5).Value = '  +4.47%  '

$ws.Cells.Item(6, 4).Value = "'0.610"
$ws.Cells.Item(6, 5).Value = '  -2.75%  '

$ws.Cells.Item(7, 5).Value = '  -0.04%  '

$ws.Cells.Item(8, 4).Value = "'55.84"
$ws.Cells.Item(8, 5).Value = '  -6.84%  '

$ws.Cells.Item(9, 4).Value = "'0.386"
$ws.Cells.Item(9, 5).Value = '  -0.42%  '

$ws.Cells.Item(10, 4).Value = "'0.0769"
$ws.Cells.Item(10, 5).Value = '  -5.20%  '

$ws.Cells.Item(11, 5).Value = '  -2.04%  '

$ws.Cells.Item(12, 4).Value = '2.303.71'
$ws.Cells.Item(12, 5).Value = '  -0.07%  '

$ws.Cells.Item(13, 4).Value = "'14.22"
$ws.Cells.Item(13, 5).Value = '  -6.45%  '

$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(14, 4).Value = "'0.800"
$ws.Cells.Item(14, 5).Value = '  -5.75%  '

$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).Value = "'20.88"
$ws.Cells.Item(15, 5).Value = '  -6.79%  '

$ws.Cells.Item(16, 5).Value = '  -4.58%  '

$ws.Cells.Item(17, 4).Value = '2.001.79'
$ws.Cells.Item(17, 5).Value = '  -0.74%  '

$ws.Cells.Item(18, 4).Value = '37.195.93'
$ws.Cells.Item(18, 5).Value = '  +0.15%  '

$ws.Cells.Item(19, 4).Value = "'69.62"
$ws.Cells.Item(19, 5).Value = '  -1.07%  '

$ws.Cells.Item(20, 5).Value = '  -3.80%  '

$ws.Cells.Item(21, 5).Value = '  -2.18%  '

$ws.Cells.Item(22, 4).Value = "'228.07"
$ws.Cells.Item(22, 5).Value = '  -1.25%  '

$ws.Cells.Item(23, 5).Value = '  +5.04%  '

$ws.Cells.Item(24, 5).Value = '  -0.04%  '

$ws.Cells.Item(25, 4).Value = "'2.34"
$ws.Cells.Item(25, 5).Value = '  -0.63%  '

$ws.Cells.Item(26, 4).Value = "'164.80"
$ws.Cells.Item(26, 5).Value = '  +0.13%  '

$ws.Cells.Item(27, 5).Value = '  -6.28%  '

$ws.Cells.Item(28, 5).Value = '  -0.73%  '

$ws.Cells.Item(29, 5).Value = '  -6.91%  '

$ws.Cells.Item(30, 4).Value = "'1.33"
$ws.Cells.Item(30, 5).Value = '  -3.61%  '

$ws.Cells.Item(31, 5).Value = '  -1.48%  '

$ws.Cells.Item(32, 5).Value = '  -4.45%  '

$ws.Cells.Item(33, 4).Value = "'0.0643"
$ws.Cells.Item(33, 5).Value = '  -2.38%  '

$ws.Cells.Item(34, 5).Value = '  +0.65%  '

$ws.Cells.Item(35, 5).Value = '  -3.88%  '

$ws.Cells.Item(36, 5).Value = '  +0.50%  '

$ws.Cells.Item(38, 5).Value = '  -4.32%  '

$ws.Cells.Item(39, 5).Value = '  -1.54%  '

$ws.Cells.Item(40, 5).Value = '  +4.07%  '

$ws.Cells.Item(41, 5).Value = '  +0.57%  '

$ws.Cells.Item(42, 5).Value = '  -5.76%  '

$ws.Cells.Item(43, 5).Value = '  -1.26%  '

$ws.Cells.Item(44, 4).Value = '1.398.79'
$ws.Cells.Item(44, 5).Value = '  +1.57%  '

$ws.Cells.Item(45, 5).Value = '  -6.11%  '

$ws.Cells.Item(46, 4).Value = "'89.30"
$ws.Cells.Item(46, 5).Value = '  -3.03%  '

$ws.Cells.Item(47, 5).Value = '  -3.74%  '

$ws.Cells.Item(48, 4).Value = "'7.02"
$ws.Cells.Item(48, 5).Value = '  -6.09%  '

$ws.Cells.Item(49, 5).Value = '  +2.06%  '

$ws.Cells.Item(50, 4).Value = '2.195.51'
$ws.Cells.Item(50, 5).Value = '  -0.11%  '

$ws.Cells.Item(51, 5).Value = '  -7.80%  '
